$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 333, shifting rows 333:346 down to 334:347.
$ws.Rows.Item(333).Insert()

# Populate the newly inserted row 333 with the new weekly price record.
$ws.Range("A333").Value = 5
$ws.Range("B333").Value = "Macroferia Regional de Talca"
$ws.Range("C333").Value = "Maule"
$ws.Range("D333").Value = 44753
$ws.Range("E333").Value = 7
$ws.Range("F333").Value = 100114014
$ws.Range("G333").Value = "Betarraga"
$ws.Range("H333").Value = "Sin especificar"
$ws.Range("I333").Value = "Primera"
$ws.Range("J333").Value = 4000
$ws.Range("K333").Value = 750
$ws.Range("L333").Value = 750
$ws.Range("M333").Value = 750
$ws.Range("N333").Value = "`$/paquete 5 unidades"
$ws.Range("O333").Value = "Región del Maule"
$ws.Range("P333").Value = 150
$ws.Range("Q333").Value = 5
$ws.Range("R333").Value = "Hortaliza"
